$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely in the edit (naive forecaster bug fix
# dropped these stale data points).
$clearRefs = @("C2","E2","C3","C4")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

# Update the remaining forecast values to the corrected, recalculated figures.
$updates = @{
    "E3" = -12.90434384910118
    "C6" = 6.277541464866965
    "E6" = 9.926356894614985
    "C7" = 4.440675708811348
    "C8" = 6.535114773304795
    "E11" = -4.72694934804071
    "E13" = -4.442195584720943
    "C15" = 1.600060471414855
    "E15" = -2.94981949413492
    "E17" = 2.83373437266663
    "E19" = 0.9898525518332146
    "C20" = 2.508469427909921
    "E20" = 2.540874511056646
    "C21" = 2.613530175870649
    "E22" = 3.025650759930021
    "E23" = -0.217407100248368
    "E24" = 6.233648892986987
    "C25" = 1.207340632843734
    "E25" = 1.091965090482994
    "E26" = 0.676128192849057
    "C28" = -0.9738659311994247
    "C29" = 3.031664248068799
    "C30" = 3.047037961814514
    "C31" = 3.770384266127391
    "C32" = 1.725130460355073
    "E32" = -1.194610791899997
    "C33" = -0.5413737111579553
    "C34" = -0.22288476972816
    "C36" = -2.221695202430862
    "E36" = -6.324519326136457
    "C37" = -0.8475469098342181
    "C38" = -1.165854108406639
    "C39" = 3.946810400327694
    "E39" = 5.401597835967098
    "C40" = -1.290816600413769
    "E40" = -8.277043919141525
    "E41" = 2.97728454211843
    "C43" = 2.93728531054509
    "E43" = 2.847433398110755
    "C44" = 2.066568225344989
    "C45" = 0.9420042531300243
    "E45" = 1.197522477793744
    "C47" = 1.648295078144368
    "E47" = 2.98747790223699
    "C48" = -1.187784794999081
    "C49" = 2.025083876303713
    "E49" = 2.460035583527898
    "C50" = 2.039329803030099
    "C51" = 0.6315208369373382
    "E51" = -0.6983252104541204
    "E52" = 8.305446157974018
    "C53" = 2.541620503438358
    "E53" = 4.63560702493917
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Host "Applied naive forecaster precision fix to $($updates.Count) cells and cleared $($clearRefs.Count) stale cells."
